$wb = $excel.ActiveWorkbook

# --- Sheet: 2018 LEAVE CREDITS ---
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Insert a new blank row above row 75 (shifts rows 75-127 down to 76-128)
$ws.Rows.Item(75).Insert()

# Expand the Table15 ListObject to cover the newly inserted row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K128"))

# Ensure the new last row's calculated formula keeps the legacy structured-reference text
$ws.Range("G128").Formula = '=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"",Table15[[#This Row],[EARNED]])'

# Fill in the data for the newly-inserted row (new leave entry)
$ws.Range("B75").Value2 = "FL(2-0-0)"
$ws.Range("D75").Value2 = 2

# Fill in the data for the row that is now row 83 (new leave entry / remark)
$ws.Range("B83").Value2 = "SP(1-0-0)"
$ws.Range("K83").Value2 = "GRAD 7/13/2023"

# Restore view state (scrolled position / selection) for this sheet
$ws.Application.Goto($ws.Range("A73"), $false)
$ws.Range("B84").Select()

# --- Sheet: 2017 LEAVE BALANCE ---
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")
$ws2.Range("A40").Value2 = 45110
$ws2.Range("B40").Select()

$wb.Application.Calculate()
